$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price updates remain stored as text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "43.677.19"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "2.202.14"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D5").Value = "257.65"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "84.35"
$ws.Range("E6").Value = "  +11.69%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").Value = "44.80"
$ws.Range("E10").Value = "  +9.69%  "
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "7.19"
$ws.Range("E12").Value = "  +4.62%  "
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "2.529.43"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "14.35"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "2.208.39"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "0.784"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "43.627.78"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "69.81"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "2.37"
$ws.Range("E22").Value = "  +9.62%  "
$ws.Range("D23").Value = "231.67"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").Value = "9.10"
$ws.Range("E24").Value = "  -3.88%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "3.54"
$ws.Range("E26").Value = "  +4.83%  "
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "39.03"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").Value = "173.67"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "20.41"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D34").Value = "5.33"
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("E37").Value = "  +5.22%  "
$ws.Range("E38").Value = "  +5.40%  "
$ws.Range("D39").Value = "12.53"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +6.12%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").Value = "63.13"
$ws.Range("E42").Value = "  +5.65%  "
$ws.Range("E43").Value = "  +5.08%  "
$ws.Range("D44").Value = "0.199"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "100.01"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0978"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").Value = "  +4.92%  "
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "0.434"
$ws.Range("E50").Value = "  -6.04%  "
$ws.Range("E51").Value = "  +8.61%  "
